$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 6 columns (cliente, data, valor, colecao, vendedor,
# regional) with 5 sample rows. The template is trimmed down to just two
# columns (cliente, valor) with a single example row.

# Drop the columns that are no longer needed: data (B), colecao (D),
# vendedor (E), regional (F) - leaving cliente (A) and valor (C, which
# becomes the new column B).
$ws.Range("F1:F6").EntireColumn.Delete()
$ws.Range("E1:E6").EntireColumn.Delete()
$ws.Range("D1:D6").EntireColumn.Delete()
$ws.Range("B1:B6").EntireColumn.Delete()

# Drop the extra sample rows, keeping only the header and one data row.
$ws.Range("A3:A6").EntireRow.Delete()

# Replace the remaining sample row with a simple example.
$ws.Range("A2").Value = "Exemplo"
$ws.Range("B2").Value = 100
